$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Replace the whole contents of the paragraph that currently contains $old
# with a freshly-built paragraph using the supplied inner XML (pPr + runs),
# so paragraph-level formatting (pStyle, spacing, indent, run formatting)
# and structural details (e.g. a leading empty <w:r/>) are preserved exactly
# instead of being collapsed by a plain Find/Replace.
function Set-ParagraphXml($old, $innerXml) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $rng = $p.Range
        if ($rng.Text.TrimEnd([char]13, [char]7) -eq $old) {
            $null = $rng.InsertXML("<w:p $wns>$innerXml</w:p>")
            return $true
        }
    }
    return $false
}

# 1) Main H1 title
$null = Set-ParagraphXml "Play Amazonia Slot for Free: Review & Features" (
    '<w:pPr><w:pStyle w:val="Heading1"/></w:pPr>' +
    '<w:r><w:t>Play Amazonia Online Slot Game for Free</w:t></w:r>'
)

# 2) "What we like" bullet list
$null = Set-ParagraphXml "Exciting mini-games for enhanced gameplay" (
    '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>' +
    '<w:r/><w:r><w:t>Classic 5x3 slot with 15 paylines</w:t></w:r>'
)

$null = Set-ParagraphXml "Unique Amazonian animal symbols" (
    '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>' +
    '<w:r/><w:r><w:t>Exciting mini-games</w:t></w:r>'
)

$null = Set-ParagraphXml "Opportunities for free spins and multiplied winnings" (
    '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>' +
    '<w:r/><w:r><w:t>Opportunity to win big with Wild and Scatter symbols</w:t></w:r>'
)

$null = Set-ParagraphXml "Varied gameplay with multiple winning combinations" (
    '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>' +
    '<w:r/><w:r><w:t>Free spins and multiplied winnings</w:t></w:r>'
)

# 3) "What we don't like" bullet list
$null = Set-ParagraphXml "Dated graphics in comparison to modern slot games" (
    '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>' +
    '<w:r/><w:r><w:t>Graphics and design appear dated</w:t></w:r>'
)

$null = Set-ParagraphXml "Clumsy 3D models for animal symbols" (
    '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>' +
    '<w:r/><w:r><w:t>Clumsy and mechanical 3D models for animal symbols</w:t></w:r>'
)

# 4) Bold title repeated near the end of the document (second occurrence)
$null = Set-ParagraphXml "Play Amazonia Slot for Free: Review & Features" (
    '<w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Amazonia Online Slot Game for Free</w:t></w:r>'
)

# 5) Italic meta-description paragraph
$null = Set-ParagraphXml "Explore the Amazon Rainforest in the Amazonia online slot game by Merkur. Enjoy mini-games, free spins, and unique symbols. Play for free now." (
    '<w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Amazonia and find out why you should play this free slot game.</w:t></w:r>'
)
